# Applies the "Updated cryptos list" GitHub Actions commit:
# refreshes Price (D) and Volume(1h) (E) columns for the crypto tracker sheet,
# and reorders a few rows (Frax/TheSandbox swap, WEMIXTOKEN/Decentraland swap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.643.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "'1.879.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'331.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'0.4717"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.93%  "
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("D9").Value = "'47.46"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.50%  "
$ws.Range("D10").Value = "'0.08022"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("D11").Value = "'1.022"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("D12").Value = "'21.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "'1.874.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("D14").Value = "'5.961"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").Value = "'7.160"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").Value = "'1.005"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "'86.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.36%  "
$ws.Range("D18").Value = "'0.00001042"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").Value = "'0.06619"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").Value = "'17.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.67%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "'27.659.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("D23").Value = "'5.497"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.30%  "
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").Value = "'2.311"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").Value = "'2.098.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.81%  "
$ws.Range("D27").Value = "'156.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.24%  "
$ws.Range("D28").Value = "'20.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").Value = "'2.087"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("D30").Value = "'5.553"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("D31").Value = "'122.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("D32").Value = "'0.9646"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").Value = "'0.09536"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").Value = "'1.455"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "'5.297"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.85%  "
$ws.Range("D37").Value = "'0.06109"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").Value = "'0.02249"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("D39").Value = "'1.228"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.89%  "
$ws.Range("D40").Value = "'8.119"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.08%  "
$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").Value = "'1.002"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.5996"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.58%  "
$ws.Range("D43").Value = "'0.1892"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").Value = "'10.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.16%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.258"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.46%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5686"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("D47").Value = "'12.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.37%  "
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("D49").Value = "'1.929"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.88%  "
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("D51").Value = "'110.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.54%  "
Write-Host "Applied cryptos list update."
